# Apply the "committed gmail compose component code" tracker update:
#  - I4: "In-Progress" -> "Done"
#  - J4: "TODO"        -> "In-Progress"
#  - view scrolled right / selection moved from H8 to J10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update the status cells in row 4 (Compose Mail row)
$ws.Range("I4").Value = "Done"
$ws.Range("J4").Value = "In-Progress"

# Update the sheet view: scroll the window so column G is the leftmost
# visible column, and move the selection to J10.
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J10").Select()
